$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.518.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.880.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "467.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +9.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.70"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +13.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.749"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -9.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.86"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.43"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.505.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.82"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.887.80"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.688.99"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.87"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.49"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.60"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +10.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.34"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +15.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "37.79"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.53"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "739.69"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +10.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.84"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.72%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.86"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.77%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.56"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.351"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.61"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +16.35%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -11.28%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.78"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.21%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.44"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.25"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.42"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.87%  "
